# Update "Result" sheet rows 2 and 5: swap the Status values and replace
# the long exception-log Comments text with short summaries (the workflow
# run no longer raises the IOException - now it correctly detects that the
# exception was/was-not thrown).
$wb = $excel.ActiveWorkbook

$wsResult = $wb.Worksheets.Item("Result")
$wsResult.Range("C2").Value = "PASS"
$wsResult.Range("E2").Value = "No exception thrown"

$wsResult.Range("C5").Value = "FAIL"
$wsResult.Range("E5").Value = "An exception was meant to be thrown"

# Rows no longer need the tall, wrapped-text height now that the Comments
# text is short - let Excel auto-fit them back to the default row height.
$wsResult.Rows.Item(2).EntireRow.AutoFit()
$wsResult.Rows.Item(5).EntireRow.AutoFit()

# Update the saved selection on the "Tests" sheet, then re-activate the
# "Result" sheet so the workbook's active tab is left exactly as it was.
$wsTests = $wb.Worksheets.Item("Tests")
$wsTests.Range("D6").Select()
$wsResult.Activate()
